$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.217.80'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '3.302.85'
$ws.Range("E3").Value = '  -2.04%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '189.54'
$ws.Range("E5").Value = '  +2.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '556.67'
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  -1.99%  '
$ws.Range("D9").Value = '3.296.55'
$ws.Range("E9").Value = '  -2.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.185'
$ws.Range("E10").Value = '  -2.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.588'
$ws.Range("E11").Value = '  -1.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.67'
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000271'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.67'
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").Value = '3.831.40'
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '616.02'
$ws.Range("E16").Value = '  +1.57%  '
$ws.Range("D17").Value = '66.175.93'
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.05'
$ws.Range("E18").Value = '  +0.39%  '
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("D20").Value = '3.297.37'
$ws.Range("E20").Value = '  -1.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.03'
$ws.Range("E21").Value = '  -5.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.909'
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.40'
$ws.Range("E23").Value = '  +8.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.63'
$ws.Range("E24").Value = '  +5.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.96'
$ws.Range("E25").Value = '  -2.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.93'
$ws.Range("E26").Value = '  -4.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.01'
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.74'
$ws.Range("E28").Value = '  -0.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.65'
$ws.Range("E29").Value = '  +1.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.68'
$ws.Range("E30").Value = '  -2.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.29'
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.10'
$ws.Range("E32").Value = '  +4.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.57'
$ws.Range("E33").Value = '  +3.27%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '565.05'
$ws.Range("E34").Value = '  +6.53%  '
$ws.Range("B35").Value = 'Cosmos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.08'
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("D37").Value = '3.804.55'
$ws.Range("E37").Value = '  -1.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '57.52'
$ws.Range("E38").Value = '  -1.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '34.23'
$ws.Range("E40").Value = '  +4.94%  '
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '0.0₃0726'
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.32'
$ws.Range("E42").Value = '  -3.74%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.73'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("E44").Value = '  +0.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.338'
$ws.Range("E45").Value = '  -4.18%  '
$ws.Range("B46").Value = 'CoreDAO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.22'
$ws.Range("E46").Value = '  -8.30%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0421'
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.25'
$ws.Range("E48").Value = '  +2.09%  '
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.57'
$ws.Range("E50").Value = '  -4.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("E51").Value = '  +0.00%  '
